# [Feat:KSW] Add Component in WeaponInfo(Bouncing)
# Adds a new "Bouncing" (int) column (N) to the WeaponInfo sheet, with
# header/type/description rows and per-item bounce-count data, and
# updates the sheet view (zoom + selection) and column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column N data -------------------------------------------------
# Row 1: field key, Row 2: field type, Row 3: Korean description
$ws.Range("N1").Value = "Bouncing"
$ws.Range("N2").Value = "int"
$ws.Range("N3").Value = "튕기는 횟수"

# Per-weapon bounce counts (rows 4-41), keyed by position.
$bounceValues = @(0,1,2,2,3, 0,1,2,2,3, 0,1,2,2,3, 0,1,2,2,3, 0,1,2,2,3, 0,1,2,2,3, 0,1,2,2,3, 1,0,0)
$row = 4
foreach ($v in $bounceValues) {
    $ws.Cells.Item($row, 14).Value = $v
    $row = $row + 1
}

# --- Formatting ----------------------------------------------------------
# N1 should carry the same bold/centered header style as the rest of row 1.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# N2/N3 should carry the same style as the rest of rows 2/3.
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# Re-apply the values (PasteSpecial(formats) should not touch values, but
# make sure nothing was clobbered).
$ws.Range("N1").Value = "Bouncing"
$ws.Range("N2").Value = "int"
$ws.Range("N3").Value = "튕기는 횟수"

# Column N width (bestfit-like width used by the rest of the sheet).
$ws.Columns.Item(14).ColumnWidth = 11

# --- Sheet view: zoom + selection ----------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("R25").Select()
